$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The 'Price' column holds plain-text numeric-looking strings in the
# source sheet (t=inlineStr); assigning a numeric-looking literal to
# .Value would otherwise get auto-coerced to a real number by Excel,
# so those specific cells are pre-formatted as Text to keep them strings.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"

$ws.Range("D2").Value = '25.298.30'
$ws.Range("E2").Value = '  -2.78%  '
$ws.Range("D3").Value = '1.567.53'
$ws.Range("E3").Value = '  -3.84%  '
$ws.Range("E4").Value = '  -0.37%  '
$ws.Range("D5").Value = '207.64'
$ws.Range("E5").Value = '  -3.02%  '
$ws.Range("E6").Value = '  -0.34%  '
$ws.Range("D8").Value = '0.243'
$ws.Range("E9").Value = '  -1.85%  '
$ws.Range("D10").Value = '17.88'
$ws.Range("E10").Value = '  -1.97%  '
$ws.Range("D11").Value = '0.0782'
$ws.Range("E11").Value = '  -0.91%  '
$ws.Range("D12").Value = '1.786.71'
$ws.Range("E12").Value = '  -3.75%  '
$ws.Range("D13").Value = '1.575.86'
$ws.Range("E13").Value = '  -3.36%  '
$ws.Range("D14").Value = '4.04'
$ws.Range("D15").Value = '0.506'
$ws.Range("E15").Value = '  -3.05%  '
$ws.Range("D16").Value = '25.296.85'
$ws.Range("E16").Value = '  -2.68%  '
$ws.Range("D17").Value = '59.62'
$ws.Range("E17").Value = '  -2.67%  '
$ws.Range("D18").Value = '0.0₃0708'
$ws.Range("E18").Value = '  -4.15%  '
$ws.Range("E19").Value = '  -0.46%  '
$ws.Range("D20").Value = '185.26'
$ws.Range("E20").Value = '  -2.30%  '
$ws.Range("E21").Value = '  -2.07%  '
$ws.Range("D22").Value = '9.31'
$ws.Range("E22").Value = '  -2.44%  '
$ws.Range("E23").Value = '  -2.81%  '
$ws.Range("E24").Value = '  -2.18%  '
$ws.Range("E25").Value = '  -0.44%  '
$ws.Range("D26").Value = '141.05'
$ws.Range("E26").Value = '  -1.89%  '
$ws.Range("D27").Value = '1.65'
$ws.Range("E27").Value = '  -7.37%  '
$ws.Range("D28").Value = '14.86'
$ws.Range("E28").Value = '  -1.66%  '
$ws.Range("D29").Value = '6.45'
$ws.Range("E29").Value = '  -4.07%  '
$ws.Range("E30").Value = '  -6.25%  '
$ws.Range("E31").Value = '  -3.79%  '
$ws.Range("E32").Value = '  -2.29%  '
$ws.Range("D33").Value = '2.99'
$ws.Range("E33").Value = '  -3.86%  '
$ws.Range("E34").Value = '  -1.48%  '
$ws.Range("D35").Value = '2.31'
$ws.Range("E35").Value = '  -4.03%  '
$ws.Range("D36").Value = '1.090.51'
$ws.Range("E36").Value = '  -3.65%  '
$ws.Range("E37").Value = '  -0.78%  '
$ws.Range("E38").Value = '  -4.75%  '
$ws.Range("E39").Value = '  -2.22%  '
$ws.Range("E40").Value = '  -3.60%  '
$ws.Range("E41").Value = '  -9.02%  '
$ws.Range("D42").Value = '0.765'
$ws.Range("E42").Value = '  -1.01%  '
$ws.Range("D43").Value = '92.42'
$ws.Range("E43").Value = '  -5.75%  '
$ws.Range("D44").Value = '5.05'
$ws.Range("E44").Value = '  -2.59%  '
$ws.Range("D45").Value = '1.702.14'
$ws.Range("E45").Value = '  -3.67%  '
$ws.Range("D46").Value = '0.0₆0112'
$ws.Range("E46").Value = '  -2.54%  '
$ws.Range("D47").Value = '52.77'
$ws.Range("E47").Value = '  -3.37%  '
$ws.Range("E48").Value = '  -4.04%  '
$ws.Range("E49").Value = '  -3.70%  '
$ws.Range("E50").Value = '  -2.05%  '
$ws.Range("E51").Value = '  -0.42%  '
